$d = $word.ActiveDocument

# Range.Text on a whole paragraph includes a trailing paragraph mark
# (CR, chr 13) and, for paragraphs inside a table cell, a trailing
# end-of-cell mark (chr 7). Strip both so label comparisons are exact.
$cr = [char]13
$bel = [char]7
$trimChars = "$cr$bel".ToCharArray()

function Clean-ParaText($s) {
    return $s.TrimEnd($trimChars)
}

# --- Add ${observations} placeholder text under the "Observations" heading ---
# The heading is followed by a run of empty paragraphs; the placeholder text
# goes into the run of the second empty paragraph, and the (now redundant)
# third empty paragraph is removed.
$obsHeading = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Clean-ParaText $d.Paragraphs.Item($i).Range.Text) -eq 'Observations') {
        $obsHeading = $i
        break
    }
}
if ($obsHeading -ne $null) {
    $target = $obsHeading + 2
    $toRemove = $obsHeading + 3
    $d.Paragraphs.Item($target).Range.Text = '${observations}'
    $d.Paragraphs.Item($toRemove).Range.Delete()
}

# --- Add ${date} placeholder text to the Validation table's Date cell ---
$dateLabel = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ((Clean-ParaText $d.Paragraphs.Item($i).Range.Text) -eq 'Date') {
        $dateLabel = $i
        break
    }
}
if ($dateLabel -ne $null) {
    $d.Paragraphs.Item($dateLabel + 1).Range.Text = '${date}'
}
